{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change being applied (commit: \"add one comment for the reflection\n// paper\"):\n//   Insert a new sentence (\"This suggests that these slip events occur\n//   where faults deform ductily in zones that are several kilometers\n//   thick and that contain substantial fluid-filled porosity. \") right\n//   before the existing sentence \"This change in reflection character\n//   may provide a new technique ... in subduction zones.\" The\n//   `_GoBack` bookmark \u2014 originally sitting alone in a later, otherwise\n//   empty paragraph \u2014 is relocated to sit mid-sentence in the new text,\n//   right after \"...deform ductily in\".\n\nconst body = context.document.body;\n\n// 1. Delete the pre-existing `_GoBack` bookmark first (it currently\n//    lives by itself in an empty paragraph further down). Doing this\n//    before inserting the replacement text avoids any ambiguity from\n//    having two bookmarks that share the same name at once.\nconst existingGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!existingGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2. Locate the exact range to replace: from the start of \"This change\n//    in reflection character...\" through the end of \"subduction\n//    zones. \" (inclusive of the trailing space). Both anchor strings\n//    are unique in the document, so this can't collide with similar\n//    phrasing elsewhere (e.g. \"...subduction thrust...\").\nconst startResults = body.search(\n  \"This change in reflection character may provide a new technique to map the landward extent of rupture in great earthquakes and improve the characterization of seismic hazards in \",\n  { matchCase: true }\n);\nstartResults.load(\"text\");\nconst endResults = body.search(\"subduction zones. \", { matchCase: true });\nendResults.load(\"text\");\nawait context.sync();\n\nif (startResults.items.length === 0 || endResults.items.length === 0) {\n  throw new Error(\"Could not locate the target sentence to edit.\");\n}\n\nconst targetRange = startResults.items[0].expandTo(endResults.items[0]);\n\n// 3. Replace that range with OOXML containing the new sentence (with\n//    the relocated `_GoBack` bookmark embedded after \"ductily in\"),\n//    followed by the original sentence, unchanged.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>\n<w:r><w:t xml:space=\"preserve\">This suggests that these slip events occur where faults deform </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/><w:r><w:t>ductily</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> in</w:t></w:r>\n<w:bookmarkStart w:id=\"100\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"100\"/>\n<w:r><w:t xml:space=\"preserve\"> zones that are se</w:t></w:r>\n<w:r><w:t>ver</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">al kilometers thick and that contain </w:t></w:r>\n<w:r><w:t>substantial</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> fluid-filled porosity. </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">This change in reflection character may provide a new technique to map the landward extent of rupture in great earthquakes and improve the characterization of seismic hazards in </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/><w:r><w:t>subduction</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> zones. </w:t></w:r>\n</w:p></w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change being applied (commit: \"add one comment for the reflection\n# paper\"):\n#   Insert a new sentence (\"This suggests that these slip events occur\n#   where faults deform ductily in zones that are several kilometers\n#   thick and that contain substantial fluid-filled porosity. \") right\n#   before the existing sentence \"This change in reflection character\n#   may provide a new technique ... in subduction zones.\" The\n#   `_GoBack` bookmark -- originally sitting alone in a later, otherwise\n#   empty paragraph -- is relocated to sit mid-sentence in the new text,\n#   right after \"...deform ductily in\".\n\n$d = $word.ActiveDocument\n\n# 1. Delete the pre-existing `_GoBack` bookmark first (it currently\n#    lives by itself in an empty paragraph further down in the\n#    document). Doing this before adding the replacement bookmark\n#    avoids any ambiguity from having two bookmarks share the same\n#    name at once.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Locate the full sentence to replace -- from \"This change in\n#    reflection character...\" through the end of \"...subduction\n#    zones. \" -- using a single unique search string, and overwrite it\n#    with the new sentence followed by the original sentence text.\n$targetRng = $d.Content\n$targetRng.Find.Execute(\"This change in reflection character may provide a new technique to map the landward extent of rupture in great earthquakes and improve the characterization of seismic hazards in subduction zones. \") | Out-Null\n\n$newText = \"This suggests that these slip events occur where faults deform ductily in zones that are several kilometers thick and that contain substantial fluid-filled porosity. This change in reflection character may provide a new technique to map the landward extent of rupture in great earthquakes and improve the characterization of seismic hazards in subduction zones. \"\n$targetRng.Text = $newText\n\n# 3. Re-insert the `_GoBack` bookmark at its new home: right after\n#    \"...deform ductily in\" and before \" zones that are several...\".\n$bmRng = $d.Content\n$bmRng.Find.Execute(\"deform ductily in\") | Out-Null\n$bmRng.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $bmRng) | Out-Null\n"}
